# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.826.21'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '2.900.53'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''199.51'
$ws.Range('E5').Value = '  +5.56%  '
$ws.Range('D6').Value = '''594.88'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.548'
$ws.Range('E8').Value = '  -1.14%  '
$ws.Range('D9').Value = '''0.197'
$ws.Range('E9').Value = '  +2.17%  '
$ws.Range('D10').Value = '2.897.70'
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('D11').Value = '''0.427'
$ws.Range('E11').Value = '  +15.52%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '''4.85'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '3.432.27'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').Value = '75.755.70'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '''27.47'
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').Value = '2.899.07'
$ws.Range('E18').Value = '  +2.83%  '
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = '''12.77'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').Value = '''371.03'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '''2.30'
$ws.Range('E22').Value = '  +1.92%  '
$ws.Range('D23').Value = '''4.26'
$ws.Range('E23').Value = '  +3.93%  '
$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '''70.90'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = '3.053.45'
$ws.Range('E26').Value = '  +2.79%  '
$ws.Range('D27').Value = '''4.17'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = '''9.60'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('E29').Value = '  +2.86%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('D32').Value = '''499.89'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = '''7.70'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '''165.32'
$ws.Range('E36').Value = '  +1.33%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('E39').Value = '  -5.57%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '''180.11'
$ws.Range('E41').Value = '  -2.55%  '
$ws.Range('D42').Value = '''0.0998'
$ws.Range('E42').Value = '  +16.04%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('E44').Value = '  -2.59%  '
$ws.Range('E45').Value = '  -2.39%  '
$ws.Range('D46').Value = '''40.13'
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('E47').Value = '  -3.30%  '
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D49').Value = '''0.568'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('E51').Value = '  +3.01%  '
